# Fix f0_voicing_dur typos: "F1" -> "f0" in the stimulus name column,
# and move the cursor/selection to B7 (as left by the author after editing).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("klatt_params")

$ws.Range("A2").Value = "Highf0_LongVoicing"
$ws.Range("A3").Value = "Highf0_ShortVoicing"
$ws.Range("A4").Value = "Lowf0_LongVoicing"
$ws.Range("A5").Value = "Lowf0_ShortVoicing"

$ws.Range("B7").Select()
